$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Rows 3-19 hold the per-training data. Column H = "PERIOD TO EXPIRE"
# (numeric, decreases by 1 day), column I = "LAST UPDATE" (text date,
# moves from 03-Nov-2025 to 04-Nov-2025).

$lastUpdateRange = $ws.Range("I3:I19")
# Force text so Excel doesn't auto-convert the date-like literal into a
# real date serial/number-format when we assign it below.
$lastUpdateRange.NumberFormat = "@"

for ($row = 3; $row -le 19; $row++) {
    $periodCell = $ws.Cells.Item($row, 8)
    $periodCell.Value2 = $periodCell.Value2 - 1

    $ws.Cells.Item($row, 9).Value2 = "04-Nov-2025"
}

# Restore the original (General) number format now that the text is safely in place.
$lastUpdateRange.NumberFormat = "General"
